$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.074.56'
$ws.Range('E2').Value = '  -2.11%  '

$ws.Range('D3').Value = '1.827.50'
$ws.Range('E3').Value = '  -0.92%  '

$ws.Range('E4').Value = '  -0.76%  '

$ws.Range('D5').Value = '''311.51'
$ws.Range('E5').Value = '  -1.78%  '

$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.78%  '

$ws.Range('D7').Value = '''0.4241'
$ws.Range('E7').Value = '  -1.24%  '

$ws.Range('D8').Value = '''0.3673'
$ws.Range('E8').Value = '  -1.66%  '

$ws.Range('D9').Value = '''0.07227'
$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').Value = '''0.8470'
$ws.Range('E10').Value = '  -2.87%  '

$ws.Range('E11').Value = '  -3.33%  '

$ws.Range('D12').Value = '1.820.94'
$ws.Range('E12').Value = '  -1.34%  '

$ws.Range('D13').Value = '''6.669'
$ws.Range('E13').Value = '  -0.60%  '

$ws.Range('E14').Value = '  -1.97%  '

$ws.Range('E15').Value = '  -1.10%  '

$ws.Range('D16').Value = '''89.58'
$ws.Range('E16').Value = '  +0.85%  '

$ws.Range('D17').Value = '''1.002'
$ws.Range('E17').Value = '  -0.95%  '

$ws.Range('D18').Value = '''0.000008761'
$ws.Range('E18').Value = '  -2.41%  '

$ws.Range('E19').Value = '  -0.70%  '

$ws.Range('D20').Value = '''14.88'
$ws.Range('E20').Value = '  -3.15%  '

$ws.Range('D21').Value = '27.128.65'
$ws.Range('E21').Value = '  -1.95%  '

$ws.Range('D22').Value = '''5.133'

$ws.Range('D23').Value = '''10.82'
$ws.Range('E23').Value = '  -1.66%  '

$ws.Range('D24').Value = '2.051.02'
$ws.Range('E24').Value = '  -0.98%  '

$ws.Range('D25').Value = '''1.980'
$ws.Range('E25').Value = '  +0.59%  '

$ws.Range('D26').Value = '''151.63'
$ws.Range('E26').Value = '  -2.01%  '

$ws.Range('D27').Value = '''2.253'
$ws.Range('E27').Value = '  +4.46%  '

$ws.Range('D28').Value = '''18.20'
$ws.Range('E28').Value = '  -1.61%  '

$ws.Range('D29').Value = '''5.252'
$ws.Range('E29').Value = '  -1.39%  '

$ws.Range('D30').Value = '''116.86'
$ws.Range('E30').Value = '  -0.72%  '

$ws.Range('D31').Value = '''0.08713'
$ws.Range('E31').Value = '  -2.15%  '

$ws.Range('D32').Value = '''1.182'
$ws.Range('E32').Value = '  -3.03%  '

$ws.Range('D33').Value = '''0.7377'
$ws.Range('E33').Value = '  -4.69%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''4.435'
$ws.Range('E34').Value = '  -1.94%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.900'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('D36').Value = '''1.000'
$ws.Range('E36').Value = '  -0.90%  '

$ws.Range('D37').Value = '''1.092'
$ws.Range('E37').Value = '  -3.28%  '

$ws.Range('D38').Value = '''0.01943'
$ws.Range('E38').Value = '  -1.54%  '

$ws.Range('E39').Value = '  -1.13%  '

$ws.Range('D40').Value = '''7.325'
$ws.Range('E40').Value = '  +2.59%  '

$ws.Range('D41').Value = '''2.875'
$ws.Range('E41').Value = '  -0.16%  '

$ws.Range('D42').Value = '''0.1688'
$ws.Range('E42').Value = '  +0.09%  '

$ws.Range('D43').Value = '''0.5076'
$ws.Range('E43').Value = '  -0.79%  '

$ws.Range('D44').Value = '''8.571'
$ws.Range('E44').Value = '  -2.13%  '

$ws.Range('D45').Value = '''1.974'
$ws.Range('E45').Value = '  +7.30%  '

$ws.Range('D46').Value = '''10.50'
$ws.Range('E46').Value = '  -1.67%  '

$ws.Range('E47').Value = '  +0.01%  '

$ws.Range('D48').Value = '''105.71'

$ws.Range('D49').Value = '''1.000'
$ws.Range('E49').Value = '  -0.92%  '

$ws.Range('D51').Value = '''1.652'
$ws.Range('E51').Value = '  -1.83%  '
